## Wed, Mar 18, 2020 11:04:29 PM
##
## 1) Slide 16 table: switch the table style applied to the cash-flow
##    summary table from the default "Table_0" style to the alternate
##    built-in style ({820F7B98-3AB4-4EC1-BD4D-0603BA50216B}).
##
## 2) Presentation theme: the deck's applied design ("Integral") is
##    swapped back out for the stock "Office Theme" palette - i.e. the
##    12 theme colours that currently paint the slides (via the Slide
##    Master's theme) are replaced with the default Office colour set.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 ------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{820F7B98-3AB4-4EC1-BD4D-0603BA50216B}")

# --- 2) Swap the applied theme's colour palette ---------------------
$theme = $p.SlideMaster.Theme
$colours = $theme.ThemeColorScheme

$officeThemeRgb = @(
    0,          # Dark 1    - 000000
    16777215,   # Light 1   - FFFFFF
    6968388,    # Dark 2    - 44546A
    15132391,   # Light 2   - E7E6E6
    13998939,   # Accent 1  - 5B9BD5
    3243501,    # Accent 2  - ED7D31
    10855845,   # Accent 3  - A5A5A5
    49407,      # Accent 4  - FFC000
    12874308,   # Accent 5  - 4472C4
    4697456,    # Accent 6  - 70AD47
    12673797,   # Hyperlink - 0563C1
    7491477     # Followed Hyperlink - 954F72
)

for ($i = 1; $i -le $colours.Count; $i++) {
    $colours.Item($i).RGB = $officeThemeRgb[$i - 1]
}
